$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '41.986.06'
$ws.Range("E2").Value = '  +5.20%  '

# Row 3
$ws.Range("D3").Value = '2.257.61'
$ws.Range("E3").Value = '  +1.84%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.33%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.11%  '

# Row 7
$ws.Range("E7").Value = '  +3.50%  '

# Row 8
$ws.Range("E8").Value = '  +0.08%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.483'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.11%  '

# Row 10
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '54.71'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.44%  '

# Row 11
$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '32.68'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.82%  '

# Row 12
$ws.Range("E12").Value = '  +2.48%  '

# Row 13
$ws.Range("E13").Value = '  +3.39%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.44%  '

# Row 15
$ws.Range("D15").Value = '2.611.61'
$ws.Range("E15").Value = '  +2.22%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.15'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.71%  '

# Row 17
$ws.Range("D17").Value = '2.261.79'
$ws.Range("E17").Value = '  +2.57%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.757'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.46%  '

# Row 19
$ws.Range("D19").Value = '41.899.02'
$ws.Range("E19").Value = '  +5.22%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +9.13%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0908'
$ws.Range("E21").Value = '  +2.41%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.61%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.24'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.54%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '241.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.82%  '

# Row 25
$ws.Range("E25").Value = '  +5.51%  '

# Row 26
$ws.Range("E26").Value = '  -0.04%  '

# Row 27
$ws.Range("E27").Value = '  +4.11%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.43%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.95%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.66'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.42%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.10'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.74%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '158.91'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.09%  '

# Row 33
$ws.Range("E33").Value = '  +0.12%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.15'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.82%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0743'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.49%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.84%  '

# Row 37
$ws.Range("E37").Value = '  +2.94%  '

# Row 38
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.105'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.52%  '

# Row 39
$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.116'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.01%  '

# Row 40
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.54'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.06%  '

# Row 41
$ws.Range("E41").Value = '  +4.31%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.94'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.33%  '

# Row 43
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.07'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +12.79%  '

# Row 44
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.051.81'
$ws.Range("E44").Value = '  -2.81%  '

# Row 45
$ws.Range("E45").Value = '  +3.77%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.12'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.26%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.90'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.89%  '

# Row 48
$ws.Range("E48").Value = '  -3.53%  '

# Row 49
$ws.Range("D49").Value = '2.483.99'
$ws.Range("E49").Value = '  +2.42%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.52'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.94%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.14'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.47%  '
